$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new column "Faithlife" with a talk date for "Async ASP.NET"
$ws.Range("M1").Value = "Faithlife"
$ws.Range("M3").Value = "2021-07"

# Move the active selection to M4, matching where the user left off editing
$ws.Range("M4").Select()
